# Noticias sobre las elecciones presidenciales.
# Populate the "JeimySosa" sheet (9th sheet) with the news items, resize
# a couple of columns, bump the row heights for the wrapped rows, and
# make this sheet the active / selected one in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(9)

# --- Row 3 -------------------------------------------------------------
$ws.Range("D3").Value = "http://www.elcolombiano.com/colombia/politica/derecha-e-izquierda-de-nuevo-las-protagonistas-FI8309447"
$ws.Range("B3").Value = "Marzo 6/2018"
$ws.Range("A3").Value = "Derecha e izquierda de nuevo las protagonistas"

# --- Row 4 (wrapped) -----------------------------------------------------
$ws.Range("A4").Value = "Santos dice que fue un ""error"" que Farc se presentaran a contienda electoral tan rápido"
$ws.Range("A4").WrapText = $true
$ws.Range("B4").Value = "Marzo 8/2018"
$ws.Range("D4").Value = "http://www.elcolombiano.com/elecciones-2018-colombia/santos-entiende-retiro-de-timochenko-de-la-contienda-electoral-FA8332061"

# --- Row 5 (wrapped) -----------------------------------------------------
$ws.Range("A5").Value = "Sin candidato presidencial, Farc se conforma con el congreso"
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "Marzo 9/2018"
$ws.Range("D5").Value = "http://www.elcolombiano.com/colombia/politica/sin-candidato-presidencial-farc-se-conforma-con-el-congreso-AB8333626"

# --- Row 6 (wrapped) -----------------------------------------------------
$ws.Range("A6").Value = "Los candidatos colombianos que copieron propagandas extrangeras"
$ws.Range("A6").WrapText = $true
$ws.Range("B6").Value = "Marzo 9/2018"
$ws.Range("D6").Value = "http://www.elcolombiano.com/elecciones-2018-colombia/candidatos-colombianos-que-copiaron-propagandas-extranjeras-FF8337878"

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "Por fin habrá tarjetón presidencial."
$ws.Range("B7").Value = "Marzo 10/2018"
$ws.Range("D7").Value = "http://www.elcolombiano.com/elecciones-2018-colombia/por-fin-habra-tarjeton-presidencial-XE8341761"

# --- Row 8 -----------------------------------------------------------------
$ws.Range("A8").Value = "El Vice, una figura con un poder creciente"
$ws.Range("B8").Value = "Marzo 10/2018"
$ws.Range("D8").Value = "http://www.elcolombiano.com/elecciones-2018-colombia/el-vice-una-figura-con-un-poder-creciente-CD8342602"

# --- Row 9 -----------------------------------------------------------------
$ws.Range("A9").Value = "Estos son los Candidatos Transparentes."
$ws.Range("B9").Value = "Marzo 10/2018"
$ws.Range("D9").Value = "http://www.elcolombiano.com/elecciones-2018-colombia/estos-son-los-candidatos-transparentes-BD8342641"

# --- Row heights for the wrapped rows ------------------------------------
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30

# --- Column widths -----------------------------------------------------
# The COM layer quantizes ColumnWidth to whole screen pixels, so these are
# the closest achievable inputs to the target stored widths of
# 42.5703125 and 13.28515625 characters.
$ws.Columns.Item(1).ColumnWidth = 41.6
$ws.Columns.Item(2).ColumnWidth = 12.45

# --- Make this the active / selected sheet --------------------------------
$ws.Activate()
[void]$ws.Range("A10").Select()
